$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 308.9
$ws.Range("I12").Value = 147.17241
$ws.Range("K12").Value = 147.17241
$ws.Range("M12").Value = 22.82758999999999
$ws.Range("H17").Value = 5265531.5
$ws.Range("J17").Value = 5265531.5
$ws.Range("L17").Value = 15796594.5
$ws.Range("N17").Value = -15796930.5
$ws.Range("H28").Value = 375.22223
$ws.Range("I28").Value = 432.85715
$ws.Range("K28").Value = 432.85715
$ws.Range("M28").Value = 52.14285000000001
$ws.Range("H69").Value = 11000
$ws.Range("J69").Value = 11000
$ws.Range("L69").Value = 33000
$ws.Range("N69").Value = -34748
$ws.Range("H70").Value = 2400
$ws.Range("J70").Value = 3000
$ws.Range("L70").Value = 9000
$ws.Range("N70").Value = -9540
$ws.Range("H72").Value = 11000
$ws.Range("J72").Value = 11000
$ws.Range("L72").Value = 99000
$ws.Range("N72").Value = -107736
$ws.Range("H73").Value = 2400
$ws.Range("J73").Value = 3000
$ws.Range("L73").Value = 9000
$ws.Range("N73").Value = -10872
$ws.Range("H112").Value = 2648.75
$ws.Range("J112").Value = 3312.5
$ws.Range("L112").Value = 9937.5
$ws.Range("N112").Value = -12153.5
$ws.Range("H125").Value = 6262947
$ws.Range("I125").Value = 48000
$ws.Range("J125").Value = 8334596
$ws.Range("K125").Value = 432000
$ws.Range("L125").Value = 75011364
$ws.Range("M125").Value = -429540
$ws.Range("N125").Value = -75016284
$ws.Range("H132").Value = 1985.625
$ws.Range("I132").Value = 1555.3928
$ws.Range("J132").Value = 4997.25
$ws.Range("K132").Value = 4666.178400000001
$ws.Range("L132").Value = 14991.75
$ws.Range("M132").Value = -2136.178400000001
$ws.Range("N132").Value = -20051.75
$ws.Range("H133").Value = 94749.5
$ws.Range("J133").Value = 94749.5
$ws.Range("L133").Value = 94749.5
$ws.Range("N133").Value = -104869.5
$ws.Range("H137").Value = 2442.3044
$ws.Range("I137").Value = 2429.4443
$ws.Range("K137").Value = 7288.3329
$ws.Range("M137").Value = -4738.3329
$ws.Range("H138").Value = 3151.2188
$ws.Range("J138").Value = 3892.5
$ws.Range("L138").Value = 11677.5
$ws.Range("N138").Value = -21957.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3303.2273
$ws.Range("I2").Value = 2824.2632
$ws.Range("K2").Value = 2824.2632
$ws.Range("M2").Value = -2711.2632
$ws.Range("H26").Value = 3002.3333
$ws.Range("I26").Value = 3002.3333
$ws.Range("K26").Value = 3002.3333
$ws.Range("M26").Value = -2672.3333
$ws.Range("H61").Value = 8622.272000000001
$ws.Range("I61").Value = 8842.380999999999
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 8842.380999999999
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -8630.380999999999
$ws.Range("N61").Value = -4424
$ws.Range("H74").Value = 5856.3335
$ws.Range("I74").Value = 4169.5557
$ws.Range("J74").Value = 10916.667
$ws.Range("K74").Value = 4169.5557
$ws.Range("L74").Value = 10916.667
$ws.Range("M74").Value = -3295.5557
$ws.Range("N74").Value = -12664.667
$ws.Range("H77").Value = 5856.3335
$ws.Range("I77").Value = 4169.5557
$ws.Range("J77").Value = 10916.667
$ws.Range("K77").Value = 20847.7785
$ws.Range("L77").Value = 54583.335
$ws.Range("M77").Value = -16479.7785
$ws.Range("N77").Value = -63319.335
$ws.Range("H88").Value = 2721.8667
$ws.Range("I88").Value = 2732
$ws.Range("J88").Value = 2710.2856
$ws.Range("K88").Value = 2732
$ws.Range("L88").Value = 2710.2856
$ws.Range("M88").Value = -2326
$ws.Range("N88").Value = -3522.2856
$ws.Range("H91").Value = 2721.8667
$ws.Range("I91").Value = 2732
$ws.Range("J91").Value = 2710.2856
$ws.Range("K91").Value = 2732
$ws.Range("L91").Value = 2710.2856
$ws.Range("M91").Value = -1328
$ws.Range("N91").Value = -5518.2856
$ws.Range("H97").Value = 1238.825
$ws.Range("I97").Value = 1246.5161
$ws.Range("J97").Value = 1212.3334
$ws.Range("K97").Value = 1246.5161
$ws.Range("L97").Value = 1212.3334
$ws.Range("M97").Value = -750.5161000000001
$ws.Range("N97").Value = -2204.3334
$ws.Range("H116").Value = 3303.2273
$ws.Range("I116").Value = 2824.2632
$ws.Range("K116").Value = 2824.2632
$ws.Range("M116").Value = -530.2631999999999
$ws.Range("H132").Value = 1594.1
$ws.Range("I132").Value = 1486.5555
$ws.Range("J132").Value = 2562
$ws.Range("K132").Value = 4459.666499999999
$ws.Range("L132").Value = 7686
$ws.Range("M132").Value = -1929.666499999999
$ws.Range("N132").Value = -12746
$ws.Range("H136").Value = 8622.272000000001
$ws.Range("I136").Value = 8842.380999999999
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 26527.143
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -23977.143
$ws.Range("N136").Value = -17100
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3303.2273
$ws.Range("I3").Value = 2824.2632
$ws.Range("K3").Value = 2824.2632
$ws.Range("M3").Value = -2710.2632
$ws.Range("H99").Value = 4467.357
$ws.Range("I99").Value = 3755.7222
$ws.Range("J99").Value = 5748.3
$ws.Range("K99").Value = 3755.7222
$ws.Range("L99").Value = 5748.3
$ws.Range("M99").Value = -2257.7222
$ws.Range("N99").Value = -8744.299999999999
$ws.Range("H132").Value = 98259.664
$ws.Range("J132").Value = 98259.664
$ws.Range("L132").Value = 98259.664
$ws.Range("N132").Value = -108379.664
$ws.Range("H134").Value = 9411.632
$ws.Range("I134").Value = 9342.412
$ws.Range("K134").Value = 28027.236
$ws.Range("M134").Value = -25492.236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H132").Value = 5139.769
$ws.Range("I132").Value = 3031
$ws.Range("K132").Value = 9093
$ws.Range("M132").Value = -6563
$ws.Range("H134").Value = 3204.9644
$ws.Range("I134").Value = 1787.5
$ws.Range("K134").Value = 5362.5
$ws.Range("M134").Value = -2827.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1019.2143
$ws.Range("J12").Value = 1574.2222
$ws.Range("L12").Value = 4722.6666
$ws.Range("N12").Value = -5068.6666
$ws.Range("H18").Value = 496.55554
$ws.Range("I18").Value = 183.625
$ws.Range("K18").Value = 550.875
$ws.Range("M18").Value = -381.875
$ws.Range("H60").Value = 525.6
$ws.Range("I60").Value = 766
$ws.Range("K60").Value = 2298
$ws.Range("M60").Value = -2047
$ws.Range("H136").Value = 900
$ws.Range("I136").Value = 900
$ws.Range("K136").Value = 2700
$ws.Range("M136").Value = 2400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11646.81
$ws.Range("I93").Value = 2036.5
$ws.Range("J93").Value = 42399.8
$ws.Range("K93").Value = 2036.5
$ws.Range("L93").Value = 42399.8
$ws.Range("M93").Value = -788.5
$ws.Range("N93").Value = -44895.8
$ws.Range("H100").Value = 4892.6665
$ws.Range("I100").Value = 4595.846
$ws.Range("J100").Value = 5375
$ws.Range("K100").Value = 4595.846
$ws.Range("L100").Value = 5375
$ws.Range("M100").Value = -4054.846
$ws.Range("N100").Value = -6457
$ws.Range("H132").Value = 13162.083
$ws.Range("I132").Value = 14994.6
$ws.Range("K132").Value = 44983.8
$ws.Range("M132").Value = -42453.8
$ws.Range("H136").Value = 6584.1885
$ws.Range("I136").Value = 6579.76
$ws.Range("J136").Value = 6658
$ws.Range("K136").Value = 19739.28
$ws.Range("L136").Value = 19974
$ws.Range("M136").Value = -17189.28
$ws.Range("N136").Value = -25074

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("H31").Value = 6999
$ws.Range("I31").Value = 6999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6651
$ws.Range("H62").Value = 7998.8
$ws.Range("I62").Value = 4997.5
$ws.Range("J62").Value = 9999.666999999999
$ws.Range("K62").Value = 4997.5
$ws.Range("L62").Value = 9999.666999999999
$ws.Range("M62").Value = -4373.5
$ws.Range("N62").Value = -11247.667
$ws.Range("H65").Value = 7998.8
$ws.Range("I65").Value = 4997.5
$ws.Range("J65").Value = 9999.666999999999
$ws.Range("K65").Value = 24987.5
$ws.Range("L65").Value = 49998.335
$ws.Range("M65").Value = -21867.5
$ws.Range("N65").Value = -56238.335
$ws.Range("H125").Value = 50712
$ws.Range("J125").Value = 50712
$ws.Range("L125").Value = 50712
$ws.Range("N125").Value = -60552
$ws.Range("H132").Value = 8658.117
$ws.Range("I132").Value = 7299.4287
$ws.Range("K132").Value = 21898.2861
$ws.Range("M132").Value = -19368.2861
$ws.Range("H136").Value = 3089.2942
$ws.Range("I136").Value = 2344.875
$ws.Range("K136").Value = 7034.625
$ws.Range("M136").Value = -4484.625
$ws.Range("M24").ClearContents()
$ws.Range("N31").ClearContents()
